# Update gh-pages to output generated at 456a3b4
# Applies the data refresh captured in the commit diff:
#  - Sheet "展览": a new event row is inserted at row 27 (pushing the
#    former rows 27-36 down to 28-37), plus small "want to go" (column F)
#    counter bumps on several unrelated rows.
#  - Sheets "演出", "本地生活", "全部类型": column F counter bumps only.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: 展览 (Exhibitions)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

# Insert a brand-new row at position 27 (existing rows 27-36 shift to 28-37)
$ws1.Rows.Item(27).Insert()

# The insert leaves A27 without the shared "index" style used by the rest
# of column A (s="1": bold, bordered, centered). Copy that formatting down
# from the cell directly below it before filling in values.
$ws1.Range("A28").Copy()
$ws1.Range("A27").PasteSpecial(-4122)

$ws1.Range("A27").Value = 26
$ws1.Range("B27").Value = "2024-05-01"
$ws1.Range("C27").Value = "上海·HATSUNE MIKU meets niko and ... 集章之旅"
$ws1.Range("D27").Value = "淮海中路775号三楼 niko and ......"
$ws1.Range("E27").Value = "2024.05.01 10:00-06.02 22:00"
$ws1.Range("F27").Value = 0
$ws1.Range("G27").Value = 46
$ws1.Range("H27").Value = "https://show.bilibili.com/platform/detail.html?id=83163"
$ws1.Range("I27").Value = "//i2.hdslb.com/bfs/openplatform/202403/9lMpza7M1711528161190.jpeg"

# Column F ("想去人数") counter bumps on rows unaffected by the insert above
$ws1.Range("F2").Value = 38
$ws1.Range("F6").Value = 1691
$ws1.Range("F9").Value = 2449
$ws1.Range("F10").Value = 706
$ws1.Range("F11").Value = 567
$ws1.Range("F13").Value = 10
$ws1.Range("F15").Value = 326
$ws1.Range("F16").Value = 204
$ws1.Range("F18").Value = 2101
$ws1.Range("F20").Value = 698

# ---------------------------------------------------------------------
# Sheet 2: 演出 (Performances)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F8").Value = 59
$ws2.Range("F11").Value = 63
$ws2.Range("F14").Value = 308
$ws2.Range("F18").Value = 152
$ws2.Range("F20").Value = 272
$ws2.Range("F24").Value = 61
$ws2.Range("F25").Value = 1760
$ws2.Range("F26").Value = 239
$ws2.Range("F27").Value = 10
$ws2.Range("F28").Value = 250

# ---------------------------------------------------------------------
# Sheet 3: 本地生活 (Local life)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 1409
$ws3.Range("F6").Value = 506
$ws3.Range("F7").Value = 172

# ---------------------------------------------------------------------
# Sheet 4: 全部类型 (All types)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1409
$ws4.Range("F5").Value = 506
$ws4.Range("F7").Value = 38
$ws4.Range("F11").Value = 1691
$ws4.Range("F13").Value = 59
$ws4.Range("F17").Value = 2449
$ws4.Range("F18").Value = 706
$ws4.Range("F19").Value = 567
$ws4.Range("F22").Value = 326
$ws4.Range("F23").Value = 63
$ws4.Range("F24").Value = 204
$ws4.Range("F27").Value = 2101
$ws4.Range("F29").Value = 698
$ws4.Range("F30").Value = 152
$ws4.Range("F33").Value = 272
$ws4.Range("F38").Value = 172
$ws4.Range("F39").Value = 1760
$ws4.Range("F41").Value = 1758
$ws4.Range("F42").Value = 239
$ws4.Range("F43").Value = 529
$ws4.Range("F47").Value = 4539
$ws4.Range("F48").Value = 108
